$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at the top of the data block (row 171/172),
# pushing the existing rows 171..251 down to 173..253.
$ws.Rows("171:172").Insert()

# New row 171: Zapallo italiano, Limache, $/caja 36 unidades, date 2021-11-10 (44510)
$ws.Range("A171").Value = 3
$ws.Range("B171").Value = "Femacal de La Calera"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").Value = 44510
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 100112032
$ws.Range("G171").Value = "Zapallo italiano"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 110
$ws.Range("K171").Value = 4000
$ws.Range("L171").Value = 4300
$ws.Range("M171").Value = 4136
$ws.Range("N171").Value = "$/caja 36 unidades"
$ws.Range("O171").Value = "Limache"
$ws.Range("P171").Value = 115
$ws.Range("Q171").Value = 36
$ws.Range("R171").Value = "Hortaliza"

# New row 172: Zapallo italiano, Región de Arica y Parinacota, $/caja 70 unidades, date 2021-11-10 (44510)
$ws.Range("A172").Value = 3
$ws.Range("B172").Value = "Femacal de La Calera"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 44510
$ws.Range("E172").Value = 5
$ws.Range("F172").Value = 100112032
$ws.Range("G172").Value = "Zapallo italiano"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 110
$ws.Range("K172").Value = 8000
$ws.Range("L172").Value = 8500
$ws.Range("M172").Value = 8273
$ws.Range("N172").Value = "$/caja 70 unidades"
$ws.Range("O172").Value = "Región de Arica y Parinacota"
$ws.Range("P172").Value = 118
$ws.Range("Q172").Value = 70
$ws.Range("R172").Value = "Hortaliza"

# Keep the D column date formatting consistent with the rest of the sheet.
$ws.Range("D171:D172").NumberFormat = "YYYY-MM-DD HH:MM:SS"
